$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-07 Saturday", "2026-02-08 Sunday"),
    @("976×7=", "318×8="),
    @("744×9=", "270×2="),
    @("830×6=", "711×5="),
    @("777×3=", "254×9="),
    @("521×7=", "887×3="),
    @("240×6=", "816×5="),
    @("722×4=", "514×8="),
    @("986×5=", "708×9="),
    @("716×4=", "842×8="),
    @("294×7=", "604×2="),
    @("933×4=", "536×9="),
    @("404×9=", "189×5="),
    @("308×7=", "957×8="),
    @("930×9=", "301×2="),
    @("953×4=", "664×2="),
    @("246×3=", "370×3="),
    @("937×4=", "288×3="),
    @("941×2=", "580×3="),
    @("965×3=", "651×8="),
    @("405×4=", "273×5="),
    @("210×5=", "794×9="),
    @("441×4=", "859×2="),
    @("136×3=", "158×2="),
    @("441×2=", "960×6="),
    @("710×5=", "939×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
